# Swarmbot IO Assignment - "Fixed the enable pin being used"
#
# The KiCad function column (G) gets reshuffled / corrected, several
# "Digital out / Digital in / Analouge in" notes in column H follow the
# corresponding G value, a stray "yes" in E33 is removed, and the pin
# that was incorrectly doubling for the Enable function (row 13 / IO27)
# has its KiCad assignment cleared with a "<-" marker added in column J
# pointing back at it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New IR sensor bank (rows 5-10) ---
$ws.Range("G5").Value  = "IR 1"
$ws.Range("H5").Value  = "Analouge in"
$ws.Range("G6").Value  = "IR 2"
$ws.Range("H6").Value  = "Analouge in"
$ws.Range("G7").Value  = "IR 3"
$ws.Range("H7").Value  = "Analouge in"
$ws.Range("G8").Value  = "IR 4"
$ws.Range("H8").Value  = "Analouge in"
$ws.Range("G9").Value  = "IR 5"
$ws.Range("H9").Value  = "Analouge in"
$ws.Range("G10").Value = "IR LEDS"
$ws.Range("H10").Value = "Digital out"

# --- Battery sense moved up ---
$ws.Range("G12").Value = "Bat I"
$ws.Range("H12").Value = "Analouge in"

# --- Pin 12 / IO27: no longer used (was wrongly tied to Enable) ---
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("J13").Value = "<-"

# --- Motor sensors / switches renamed & reshuffled ---
$ws.Range("G14").Value = "Right Motor Sensor"
$ws.Range("G15").Value = "5V Switch 1"
$ws.Range("H15").Value = "Digital out"
$ws.Range("G17").Value = "Line Sensor 5"

# --- SDIO pins now carry Motor PWM / Line sensor / switch duties ---
$ws.Range("G18").Value = "Motor PWM 2"
$ws.Range("H18").Value = "Digital out (PWM)"
$ws.Range("G19").Value = "Line Sensor LED"
$ws.Range("H19").Value = "Digital out"
$ws.Range("G20").Value = "Motor PWM 1"
$ws.Range("H20").Value = "Digital out (PWM)"
$ws.Range("G21").Value = "Left Switch"
$ws.Range("H21").Value = "Digital in"
$ws.Range("G22").Value = "Motor PWM 3"
$ws.Range("H22").Value = "Digital out (PWM)"
$ws.Range("G23").Value = "Line Sensor 4"
$ws.Range("H23").Value = "Analouge in"

# --- Remaining line sensors / sonar / PCB versions reassigned ---
$ws.Range("G24").Value = "Motor PWM 4"
$ws.Range("H24").Value = "Digital out (PWM)"
$ws.Range("G25").Value = "Sonar Echo"
$ws.Range("H25").Value = "Digital in"
$ws.Range("G27").Value = "Sonar Trigger"
$ws.Range("G28").Value = "Line Sensor 3"
$ws.Range("G29").Value = "Line Sensor 2"
$ws.Range("H29").Value = "Analouge in"
$ws.Range("G30").Value = "Bottom PCB Version"
$ws.Range("H30").Value = "Analouge in"
$ws.Range("G31").Value = "Right Switch"
$ws.Range("G32").Value = "Line Sensor 1"

# --- Stray PWM flag removed from the n/c row ---
$ws.Range("E33").ClearContents()

$ws.Range("G34").Value = "Left Motor Sensor"
$ws.Range("G37").Value = "Pen Servo"
$ws.Range("G38").Value = "Top PCB Version"

# Reselect the cell the author was pointing at when they made the fix.
$ws.Range("J13").Select()
